# Edit: insert a new weekly price record as row 40 in the "Papa" sheet,
# shifting all subsequent rows (old 40..142) down by one (to 41..143).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 40; this shifts old rows 40-142 to 41-143
# and also updates the used range / dimension (A1:R142 -> A1:R143).
$ws.Rows("40:40").Insert()

# Populate the newly inserted row 40 with the new record's data.
$ws.Cells.Item(40, 1).Value  = 1
$ws.Cells.Item(40, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(40, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(40, 4).Value  = (Get-Date -Year 2023 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(40, 5).Value  = 15
$ws.Cells.Item(40, 6).Value  = 100114001
$ws.Cells.Item(40, 7).Value  = "Papa"
$ws.Cells.Item(40, 8).Value  = "Asterix"
$ws.Cells.Item(40, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(40, 10).Value = 1300
$ws.Cells.Item(40, 11).Value = 14000
$ws.Cells.Item(40, 12).Value = 15000
$ws.Cells.Item(40, 13).Value = 14462
$ws.Cells.Item(40, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(40, 15).Value = "Región Metropolitana"
$ws.Cells.Item(40, 16).Value = 578
$ws.Cells.Item(40, 17).Value = 25
$ws.Cells.Item(40, 18).Value = "Hortaliza"
